$d = $word.ActiveDocument

$replacements = @(
    @{old="76×15="; new="52×93="},
    @{old="31×62="; new="24×65="},
    @{old="89×15="; new="66×47="},
    @{old="46×92="; new="93×84="},
    @{old="39×48="; new="88×97="},
    @{old="28×48="; new="61×53="},
    @{old="58×18="; new="71×12="},
    @{old="82×71="; new="62×21="},
    @{old="35×20="; new="96×87="},
    @{old="12×13="; new="51×45="},
    @{old="62×98="; new="93×68="},
    @{old="56×72="; new="27×76="},
    @{old="75×23="; new="56×75="},
    @{old="96×55="; new="83×37="},
    @{old="44×37="; new="42×60="},
    @{old="29×23="; new="47×65="},
    @{old="27×89="; new="51×46="},
    @{old="67×91="; new="11×70="},
    @{old="68×41="; new="67×45="},
    @{old="39×57="; new="56×33="},
    @{old="43×88="; new="64×56="},
    @{old="79×15="; new="57×77="},
    @{old="77×57="; new="27×75="},
    @{old="97×41="; new="68×59="},
    @{old="37×13="; new="65×28="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
